$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 11363757
$ws.Range("I6").Value = 17857292
$ws.Range("K6").Value = 53571876
$ws.Range("M6").Value = -53571764
$ws.Range("H11").Value = 176.625
$ws.Range("I11").Value = 176.625
$ws.Range("K11").Value = 176.625
$ws.Range("M11").Value = -36.625
$ws.Range("H15").Value = 1471.7819
$ws.Range("I15").Value = 1471.7819
$ws.Range("K15").Value = 4415.3457
$ws.Range("M15").Value = -4246.3457
$ws.Range("H17").Value = 54097.42
$ws.Range("J17").Value = 57036.723
$ws.Range("L17").Value = 171110.169
$ws.Range("N17").Value = -171446.169
$ws.Range("H18").Value = 2024.2858
$ws.Range("I18").Value = 1248
$ws.Range("K18").Value = 1248
$ws.Range("M18").Value = -964
$ws.Range("H19").Value = 5001090.5
$ws.Range("I19").Value = 883.3333
$ws.Range("J19").Value = 7144036.5
$ws.Range("K19").Value = 883.3333
$ws.Range("L19").Value = 7144036.5
$ws.Range("M19").Value = -708.3333
$ws.Range("N19").Value = -7144386.5
$ws.Range("H80").Value = 3057501.2
$ws.Range("I80").Value = 1508942.1
$ws.Range("J80").Value = 5294308.5
$ws.Range("K80").Value = 4526826.300000001
$ws.Range("L80").Value = 15882925.5
$ws.Range("M80").Value = -4525828.300000001
$ws.Range("N80").Value = -15884921.5
$ws.Range("H83").Value = 3057501.2
$ws.Range("I83").Value = 1508942.1
$ws.Range("J83").Value = 5294308.5
$ws.Range("K83").Value = 13580478.9
$ws.Range("L83").Value = 47648776.5
$ws.Range("M83").Value = -13575486.9
$ws.Range("N83").Value = -47658760.5
$ws.Range("H86").Value = 5406.5713
$ws.Range("I86").Value = 4410.1113
$ws.Range("K86").Value = 4410.1113
$ws.Range("M86").Value = -3287.1113
$ws.Range("H89").Value = 5406.5713
$ws.Range("I89").Value = 4410.1113
$ws.Range("K89").Value = 22050.5565
$ws.Range("M89").Value = -16434.5565
$ws.Range("H106").Value = 5819.278
$ws.Range("I106").Value = 6754.9287
$ws.Range("J106").Value = 2544.5
$ws.Range("K106").Value = 6754.9287
$ws.Range("L106").Value = 2544.5
$ws.Range("M106").Value = -6123.9287
$ws.Range("N106").Value = -3806.5
$ws.Range("H111").Value = 3321.1
$ws.Range("I111").Value = 2816.1428
$ws.Range("J111").Value = 4499.3335
$ws.Range("K111").Value = 8448.428400000001
$ws.Range("L111").Value = 13498.0005
$ws.Range("M111").Value = -5381.428400000001
$ws.Range("N111").Value = -19632.0005
$ws.Range("H132").Value = 1715.4359
$ws.Range("I132").Value = 1612.8
$ws.Range("J132").Value = 2613.5
$ws.Range("K132").Value = 4838.4
$ws.Range("L132").Value = 7840.5
$ws.Range("M132").Value = -2308.4
$ws.Range("N132").Value = -12900.5
$ws.Range("H137").Value = 13160143
$ws.Range("I137").Value = 19232594
$ws.Range("J137").Value = 3165.5
$ws.Range("K137").Value = 57697782
$ws.Range("L137").Value = 9496.5
$ws.Range("M137").Value = -57695232
$ws.Range("N137").Value = -14596.5
$ws.Range("H138").Value = 3421.4854
$ws.Range("I138").Value = 1564.9434
$ws.Range("J138").Value = 5389.42
$ws.Range("K138").Value = 4694.8302
$ws.Range("L138").Value = 16168.26
$ws.Range("M138").Value = 445.1697999999997
$ws.Range("N138").Value = -26448.26
$ws.Range("H141").Value = 2261.39
$ws.Range("I141").Value = 1701.6072
$ws.Range("J141").Value = 12710.667
$ws.Range("K141").Value = 5104.821599999999
$ws.Range("L141").Value = 38132.001
$ws.Range("M141").Value = 75.17840000000069
$ws.Range("N141").Value = -48492.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2063.5186
$ws.Range("I2").Value = 1943.3529
$ws.Range("J2").Value = 2267.8
$ws.Range("K2").Value = 1943.3529
$ws.Range("L2").Value = 2267.8
$ws.Range("M2").Value = -1830.3529
$ws.Range("N2").Value = -2493.8
$ws.Range("H32").Value = 6033.857
$ws.Range("I32").Value = 5480.271
$ws.Range("K32").Value = 5480.271
$ws.Range("M32").Value = -5193.271
$ws.Range("H45").Value = 2804369
$ws.Range("I45").Value = 4763721
$ws.Range("K45").Value = 4763721
$ws.Range("M45").Value = -4763344
$ws.Range("H61").Value = 17076538
$ws.Range("I61").Value = 21215342
$ws.Range("K61").Value = 21215342
$ws.Range("M61").Value = -21215130
$ws.Range("H62").Value = 99989
$ws.Range("J62").Value = 99989
$ws.Range("L62").Value = 99989
$ws.Range("N62").Value = -101237
$ws.Range("H65").Value = 99989
$ws.Range("J65").Value = 99989
$ws.Range("L65").Value = 299967
$ws.Range("N65").Value = -306207
$ws.Range("H74").Value = 1954.2759
$ws.Range("I74").Value = 1920.9048
$ws.Range("J74").Value = 2041.875
$ws.Range("K74").Value = 1920.9048
$ws.Range("L74").Value = 2041.875
$ws.Range("M74").Value = -1046.9048
$ws.Range("N74").Value = -3789.875
$ws.Range("H77").Value = 1954.2759
$ws.Range("I77").Value = 1920.9048
$ws.Range("J77").Value = 2041.875
$ws.Range("K77").Value = 9604.523999999999
$ws.Range("L77").Value = 10209.375
$ws.Range("M77").Value = -5236.523999999999
$ws.Range("N77").Value = -18945.375
$ws.Range("H81").Value = 94999.5
$ws.Range("J81").Value = 94999.5
$ws.Range("L81").Value = 94999.5
$ws.Range("N81").Value = -96995.5
$ws.Range("H84").Value = 94999.5
$ws.Range("J84").Value = 94999.5
$ws.Range("L84").Value = 284998.5
$ws.Range("N84").Value = -294982.5
$ws.Range("H96").Value = 31847
$ws.Range("J96").Value = 31847
$ws.Range("L96").Value = 31847
$ws.Range("N96").Value = -37339
$ws.Range("H97").Value = 2444.4443
$ws.Range("I97").Value = 2067.8
$ws.Range("K97").Value = 2067.8
$ws.Range("M97").Value = -1571.8
$ws.Range("H102").Value = 3322.611
$ws.Range("I102").Value = 1982.7
$ws.Range("K102").Value = 1982.7
$ws.Range("M102").Value = -360.7
$ws.Range("H110").Value = 4036.9285
$ws.Range("I110").Value = 4320.737
$ws.Range("K110").Value = 4320.737
$ws.Range("M110").Value = -2275.737
$ws.Range("H116").Value = 2063.5186
$ws.Range("I116").Value = 1943.3529
$ws.Range("J116").Value = 2267.8
$ws.Range("K116").Value = 1943.3529
$ws.Range("L116").Value = 2267.8
$ws.Range("M116").Value = 350.6470999999999
$ws.Range("N116").Value = -6855.8
$ws.Range("H122").Value = 3705.5625
$ws.Range("I122").Value = 2949.3572
$ws.Range("K122").Value = 8848.071599999999
$ws.Range("M122").Value = -6398.071599999999
$ws.Range("H132").Value = 3567.625
$ws.Range("I132").Value = 3445.7222
$ws.Range("J132").Value = 3933.3333
$ws.Range("K132").Value = 10337.1666
$ws.Range("L132").Value = 11799.9999
$ws.Range("M132").Value = -7807.1666
$ws.Range("N132").Value = -16859.9999
$ws.Range("H136").Value = 17076538
$ws.Range("I136").Value = 21215342
$ws.Range("K136").Value = 63646026
$ws.Range("M136").Value = -63643476

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2063.5186
$ws.Range("I3").Value = 1943.3529
$ws.Range("J3").Value = 2267.8
$ws.Range("K3").Value = 1943.3529
$ws.Range("L3").Value = 2267.8
$ws.Range("M3").Value = -1829.3529
$ws.Range("N3").Value = -2495.8
$ws.Range("H60").Value = 49666.332
$ws.Range("J60").Value = 49666.332
$ws.Range("L60").Value = 49666.332
$ws.Range("N60").Value = -50864.332
$ws.Range("H86").Value = 24433.936
$ws.Range("J86").Value = 5161.3335
$ws.Range("L86").Value = 5161.3335
$ws.Range("N86").Value = -7407.3335
$ws.Range("H89").Value = 24433.936
$ws.Range("J89").Value = 5161.3335
$ws.Range("L89").Value = 25806.6675
$ws.Range("N89").Value = -37038.6675
$ws.Range("H105").Value = 850983.9399999999
$ws.Range("I105").Value = 1347628.5
$ws.Range("K105").Value = 1347628.5
$ws.Range("M105").Value = -1345881.5
$ws.Range("H107").Value = 2724.6785
$ws.Range("I107").Value = 3096.7222
$ws.Range("K107").Value = 3096.7222
$ws.Range("M107").Value = -1176.7222
$ws.Range("H134").Value = 1794.8572
$ws.Range("I134").Value = 1877.6111
$ws.Range("J134").Value = 1298.3334
$ws.Range("K134").Value = 5632.8333
$ws.Range("L134").Value = 3895.0002
$ws.Range("M134").Value = -3097.8333
$ws.Range("N134").Value = -8965.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22730918
$ws.Range("I31").Value = 45457172
$ws.Range("J31").Value = 4661.409
$ws.Range("K31").Value = 45457172
$ws.Range("L31").Value = 4661.409
$ws.Range("M31").Value = -45456877
$ws.Range("N31").Value = -5251.409
$ws.Range("H34").Value = 22730918
$ws.Range("I34").Value = 45457172
$ws.Range("J34").Value = 4661.409
$ws.Range("K34").Value = 45457172
$ws.Range("L34").Value = 4661.409
$ws.Range("M34").Value = -45456970
$ws.Range("N34").Value = -5065.409
$ws.Range("H58").Value = 2219.147
$ws.Range("I58").Value = 2097.8333
$ws.Range("J58").Value = 2510.3
$ws.Range("K58").Value = 2097.8333
$ws.Range("L58").Value = 2510.3
$ws.Range("M58").Value = -1894.8333
$ws.Range("N58").Value = -2916.3
$ws.Range("H62").Value = 12351348
$ws.Range("I62").Value = 2566.6
$ws.Range("K62").Value = 2566.6
$ws.Range("M62").Value = -1942.6
$ws.Range("H65").Value = 12351348
$ws.Range("I65").Value = 2566.6
$ws.Range("K65").Value = 12833
$ws.Range("M65").Value = -9713
$ws.Range("H68").Value = 63212.43
$ws.Range("J68").Value = 63747.832
$ws.Range("L68").Value = 63747.832
$ws.Range("N68").Value = -65245.832
$ws.Range("H71").Value = 63212.43
$ws.Range("J71").Value = 63747.832
$ws.Range("L71").Value = 191243.496
$ws.Range("N71").Value = -198731.496
$ws.Range("H74").Value = 67666.336
$ws.Range("J74").Value = 67666.336
$ws.Range("L74").Value = 67666.336
$ws.Range("N74").Value = -69414.336
$ws.Range("H77").Value = 67666.336
$ws.Range("J77").Value = 67666.336
$ws.Range("L77").Value = 202999.008
$ws.Range("N77").Value = -211735.008
$ws.Range("H86").Value = 10031.3125
$ws.Range("I86").Value = 11305.091
$ws.Range("K86").Value = 11305.091
$ws.Range("M86").Value = -10182.091
$ws.Range("H89").Value = 10031.3125
$ws.Range("I89").Value = 11305.091
$ws.Range("K89").Value = 56525.455
$ws.Range("M89").Value = -50909.455
$ws.Range("H99").Value = 21648.072
$ws.Range("I99").Value = 11370.637
$ws.Range("K99").Value = 11370.637
$ws.Range("M99").Value = -9872.637000000001
$ws.Range("H117").Value = 99999
$ws.Range("J117").Value = 99999
$ws.Range("L117").Value = 99999
$ws.Range("N117").Value = -109177
$ws.Range("H126").Value = 21648.072
$ws.Range("I126").Value = 11370.637
$ws.Range("K126").Value = 34111.911
$ws.Range("M126").Value = -31641.911
$ws.Range("H132").Value = 2031.8889
$ws.Range("I132").Value = 2009.3125
$ws.Range("J132").Value = 2212.5
$ws.Range("K132").Value = 6027.9375
$ws.Range("L132").Value = 6637.5
$ws.Range("M132").Value = -3497.9375
$ws.Range("N132").Value = -11697.5
$ws.Range("H134").Value = 1653.5483
$ws.Range("I134").Value = 1483.28
$ws.Range("J134").Value = 2363
$ws.Range("K134").Value = 4449.84
$ws.Range("L134").Value = 7089
$ws.Range("M134").Value = -1914.84
$ws.Range("N134").Value = -12159
$ws.Range("H136").Value = 2219.147
$ws.Range("I136").Value = 2097.8333
$ws.Range("J136").Value = 2510.3
$ws.Range("K136").Value = 6293.499899999999
$ws.Range("L136").Value = 7530.900000000001
$ws.Range("M136").Value = -3743.499899999999
$ws.Range("N136").Value = -12630.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1428807
$ws.Range("I4").Value = 274.83334
$ws.Range("K4").Value = 824.5000200000001
$ws.Range("M4").Value = -712.5000200000001
$ws.Range("H5").Value = 559.9167
$ws.Range("I5").Value = 840.8
$ws.Range("J5").Value = 359.2857
$ws.Range("K5").Value = 2522.4
$ws.Range("L5").Value = 1077.8571
$ws.Range("M5").Value = -2410.4
$ws.Range("N5").Value = -1301.8571
$ws.Range("H7").Value = 226.21428
$ws.Range("J7").Value = 127.77778
$ws.Range("L7").Value = 383.33334
$ws.Range("N7").Value = -607.33334
$ws.Range("H12").Value = 2233.2083
$ws.Range("J12").Value = 2942.7856
$ws.Range("L12").Value = 8828.356800000001
$ws.Range("N12").Value = -9174.356800000001
$ws.Range("H29").Value = 8956.571
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 10366
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 31098
$ws.Range("N29").Value = -31652
$ws.Range("M29").Value = -1223
$ws.Range("H49").Value = 8989.25
$ws.Range("I49").Value = 870
$ws.Range("J49").Value = 22521.334
$ws.Range("K49").Value = 2610
$ws.Range("L49").Value = 67564.00199999999
$ws.Range("N49").Value = -67876.00199999999
$ws.Range("M49").Value = -2454
$ws.Range("H58").Value = 10691
$ws.Range("I58").Value = 6666.3335
$ws.Range("K58").Value = 19999.0005
$ws.Range("M58").Value = -19871.0005
$ws.Range("H70").Value = 15577.637
$ws.Range("J70").Value = 18420
$ws.Range("L70").Value = 55260
$ws.Range("N70").Value = -55890
$ws.Range("H73").Value = 15577.637
$ws.Range("J73").Value = 18420
$ws.Range("L73").Value = 55260
$ws.Range("N73").Value = -57444
$ws.Range("H80").Value = 27781444
$ws.Range("J80").Value = 5332.6665
$ws.Range("L80").Value = 15997.9995
$ws.Range("N80").Value = -17869.9995
$ws.Range("H83").Value = 27781444
$ws.Range("J83").Value = 5332.6665
$ws.Range("L83").Value = 47993.9985
$ws.Range("N83").Value = -57353.9985
$ws.Range("H92").Value = 48.555557
$ws.Range("I92").Value = 20.5
$ws.Range("J92").Value = 56.57143
$ws.Range("K92").Value = 61.5
$ws.Range("L92").Value = 169.71429
$ws.Range("M92").Value = 1186.5
$ws.Range("N92").Value = -2665.71429
$ws.Range("H93").Value = 11078.723
$ws.Range("J93").Value = 13021.134
$ws.Range("L93").Value = 39063.402
$ws.Range("N93").Value = -42807.402
$ws.Range("H108").Value = 15582.75
$ws.Range("I108").Value = 9666
$ws.Range("K108").Value = 28998
$ws.Range("M108").Value = -26118
$ws.Range("H109").Value = 7730.304
$ws.Range("I109").Value = 591.8125
$ws.Range("K109").Value = 1775.4375
$ws.Range("M109").Value = -735.4375
$ws.Range("H122").Value = 610
$ws.Range("I122").Value = 512.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4612.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2162.5
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 15904.286
$ws.Range("I125").Value = 15999
$ws.Range("J125").Value = 15833.25
$ws.Range("K125").Value = 47997
$ws.Range("L125").Value = 47499.75
$ws.Range("M125").Value = -43077
$ws.Range("N125").Value = -57339.75
$ws.Range("H132").Value = 2105.0557
$ws.Range("I132").Value = 1779.6
$ws.Range("K132").Value = 16016.4
$ws.Range("M132").Value = -13486.4
$ws.Range("H135").Value = 559.9167
$ws.Range("I135").Value = 840.8
$ws.Range("J135").Value = 359.2857
$ws.Range("K135").Value = 7567.2
$ws.Range("L135").Value = 3233.5713
$ws.Range("M135").Value = -5032.2
$ws.Range("N135").Value = -8303.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 7121.3335
$ws.Range("I31").Value = 10031
$ws.Range("J31").Value = 5666.5
$ws.Range("K31").Value = 10031
$ws.Range("L31").Value = 5666.5
$ws.Range("M31").Value = -9739
$ws.Range("N31").Value = -6250.5
$ws.Range("H35").Value = 67221.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 67221.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 67221.5
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -67817.5
$ws.Range("H37").Value = 7121.3335
$ws.Range("I37").Value = 10031
$ws.Range("J37").Value = 5666.5
$ws.Range("K37").Value = 10031
$ws.Range("L37").Value = 5666.5
$ws.Range("M37").Value = -9754
$ws.Range("N37").Value = -6220.5
$ws.Range("H70").Value = 8610.964
$ws.Range("I70").Value = 7874
$ws.Range("K70").Value = 7874
$ws.Range("M70").Value = -7604
$ws.Range("H73").Value = 8610.964
$ws.Range("I73").Value = 7874
$ws.Range("K73").Value = 7874
$ws.Range("M73").Value = -6938
$ws.Range("H80").Value = 66668936
$ws.Range("I80").Value = 1750
$ws.Range("J80").Value = 120002690
$ws.Range("K80").Value = 1750
$ws.Range("L80").Value = 120002690
$ws.Range("M80").Value = -752
$ws.Range("N80").Value = -120004686
$ws.Range("H83").Value = 66668936
$ws.Range("I83").Value = 1750
$ws.Range("J83").Value = 120002690
$ws.Range("K83").Value = 8750
$ws.Range("L83").Value = 600013450
$ws.Range("M83").Value = -3758
$ws.Range("N83").Value = -600023434
$ws.Range("H97").Value = 744.5925999999999
$ws.Range("J97").Value = 1470
$ws.Range("L97").Value = 1470
$ws.Range("N97").Value = -2462
$ws.Range("H102").Value = 4230.5
$ws.Range("I102").Value = 4230.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4230.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2608.5
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 1030770.75
$ws.Range("I113").Value = 2498
$ws.Range("J113").Value = 1685126.1
$ws.Range("K113").Value = 2498
$ws.Range("L113").Value = 1685126.1
$ws.Range("M113").Value = -328
$ws.Range("N113").Value = -1689466.1
$ws.Range("H118").Value = 99999
$ws.Range("J118").Value = 99999
$ws.Range("L118").Value = 99999
$ws.Range("N118").Value = -103313
$ws.Range("H122").Value = 5355.5835
$ws.Range("I122").Value = 5030.7
$ws.Range("K122").Value = 15092.1
$ws.Range("M122").Value = -12642.1
$ws.Range("H126").Value = 2895.0527
$ws.Range("I126").Value = 2127.0715
$ws.Range("J126").Value = 5045.4
$ws.Range("K126").Value = 6381.2145
$ws.Range("L126").Value = 15136.2
$ws.Range("M126").Value = -3911.2145
$ws.Range("N126").Value = -20076.2
$ws.Range("H132").Value = 2816.7354
$ws.Range("I132").Value = 2861.2856
$ws.Range("J132").Value = 2608.8333
$ws.Range("K132").Value = 8583.856800000001
$ws.Range("L132").Value = 7826.499899999999
$ws.Range("M132").Value = -6053.856800000001
$ws.Range("N132").Value = -12886.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5530.6665
$ws.Range("I7").Value = 2996.6365
$ws.Range("K7").Value = 2996.6365
$ws.Range("M7").Value = -2884.6365
$ws.Range("H22").Value = 3116.7273
$ws.Range("I22").Value = 2549
$ws.Range("J22").Value = 3798
$ws.Range("K22").Value = 2549
$ws.Range("L22").Value = 3798
$ws.Range("M22").Value = -2254
$ws.Range("N22").Value = -4388
$ws.Range("H27").Value = 3116.7273
$ws.Range("I27").Value = 2549
$ws.Range("J27").Value = 3798
$ws.Range("K27").Value = 2549
$ws.Range("L27").Value = 3798
$ws.Range("M27").Value = -2442
$ws.Range("N27").Value = -4012
$ws.Range("H40").Value = 3404.6584
$ws.Range("I40").Value = 3390.8823
$ws.Range("K40").Value = 3390.8823
$ws.Range("M40").Value = -3254.8823
$ws.Range("H55").Value = 674.5769
$ws.Range("I55").Value = 406.45456
$ws.Range("J55").Value = 871.2
$ws.Range("K55").Value = 406.45456
$ws.Range("L55").Value = 871.2
$ws.Range("M55").Value = -233.45456
$ws.Range("N55").Value = -1217.2
$ws.Range("H92").Value = 56000
$ws.Range("I92").Value = 56000
$ws.Range("K92").Value = 56000
$ws.Range("M92").Value = -53504
$ws.Range("H126").Value = 5530.6665
$ws.Range("I126").Value = 2996.6365
$ws.Range("K126").Value = 8989.9095
$ws.Range("M126").Value = -6519.9095
$ws.Range("H132").Value = 3453.641
$ws.Range("I132").Value = 2044.6538
$ws.Range("K132").Value = 6133.9614
$ws.Range("M132").Value = -3603.9614
$ws.Range("H136").Value = 6062.3335
$ws.Range("I136").Value = 2496
$ws.Range("J136").Value = 9628.666999999999
$ws.Range("K136").Value = 7488
$ws.Range("L136").Value = 28886.001
$ws.Range("M136").Value = -4938
$ws.Range("N136").Value = -33986.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2399.5312
$ws.Range("I107").Value = 1944.091
$ws.Range("J107").Value = 3401.5
$ws.Range("K107").Value = 5832.272999999999
$ws.Range("L107").Value = 10204.5
$ws.Range("M107").Value = -3912.272999999999
$ws.Range("N107").Value = -14044.5
$ws.Range("H126").Value = 2186.9033
$ws.Range("I126").Value = 1833.3462
$ws.Range("K126").Value = 5500.0386
$ws.Range("M126").Value = -3030.0386
$ws.Range("H132").Value = 1141.1818
$ws.Range("I132").Value = 1175.8064
$ws.Range("K132").Value = 3527.4192
$ws.Range("M132").Value = -997.4191999999998
$ws.Range("H136").Value = 1130.875
$ws.Range("I136").Value = 625.6316
$ws.Range("K136").Value = 1876.8948
$ws.Range("M136").Value = 673.1052
